# Scheduled-runner refresh of market-price / profit columns (H..N) across
# the leve-profit sheets (ALC, ARM, CRP, CUL, GSM, LTW, WVR). Values are
# plain scraped numbers (no formulas in this workbook), so each changed
# cell is just re-written with its new value; a couple of rows gain or
# lose a LeveProfit (M/N) cell entirely, handled with ClearContents().
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H9").Value = 100.77778
$ws_ALC.Range("I9").Value = 120
$ws_ALC.Range("J9").Value = 76.75
$ws_ALC.Range("K9").Value = 120
$ws_ALC.Range("L9").Value = 76.75
$ws_ALC.Range("M9").Value = 49
$ws_ALC.Range("N9").Value = -414.75
$ws_ALC.Range("H12").Value = 175
$ws_ALC.Range("I12").Value = 175
$ws_ALC.Range("K12").Value = 175
$ws_ALC.Range("M12").Value = -5
$ws_ALC.Range("H28").Value = 475
$ws_ALC.Range("I28").Value = 200
$ws_ALC.Range("J28").Value = 566.6667
$ws_ALC.Range("K28").Value = 200
$ws_ALC.Range("L28").Value = 566.6667
$ws_ALC.Range("M28").Value = 285
$ws_ALC.Range("N28").Value = -1536.6667
$ws_ALC.Range("H33").Value = 369.9
$ws_ALC.Range("I33").Value = 310
$ws_ALC.Range("J33").Value = 709.3333
$ws_ALC.Range("K33").Value = 310
$ws_ALC.Range("L33").Value = 709.3333
$ws_ALC.Range("M33").Value = -81
$ws_ALC.Range("N33").Value = -1167.3333
$ws_ALC.Range("H40").Value = 4405.3335
$ws_ALC.Range("J40").Value = 5497
$ws_ALC.Range("L40").Value = 5497
$ws_ALC.Range("N40").Value = -5847
$ws_ALC.Range("H64").Value = 5250
$ws_ALC.Range("J64").Value = 5250
$ws_ALC.Range("L64").Value = 5250
$ws_ALC.Range("N64").Value = -5746
$ws_ALC.Range("H67").Value = 5250
$ws_ALC.Range("J67").Value = 5250
$ws_ALC.Range("L67").Value = 5250
$ws_ALC.Range("N67").Value = -6966
$ws_ALC.Range("H70").Value = 4166.6924
$ws_ALC.Range("J70").Value = 4097.3335
$ws_ALC.Range("L70").Value = 12292.0005
$ws_ALC.Range("N70").Value = -12832.0005
$ws_ALC.Range("H73").Value = 4166.6924
$ws_ALC.Range("J73").Value = 4097.3335
$ws_ALC.Range("L73").Value = 12292.0005
$ws_ALC.Range("N73").Value = -14164.0005
$ws_ALC.Range("H100").Value = 2448.9
$ws_ALC.Range("I100").Value = 2569.8572
$ws_ALC.Range("K100").Value = 2569.8572
$ws_ALC.Range("M100").Value = -2028.8572
$ws_ALC.Range("H116").Value = 5998
$ws_ALC.Range("I116").Value = 5998
$ws_ALC.Range("K116").Value = 5998
$ws_ALC.Range("M116").Value = -2556

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H74").Value = 753.75
$ws_ARM.Range("J74").Value = 0
$ws_ARM.Range("L74").Value = 0
$ws_ARM.Range("N74").ClearContents()
$ws_ARM.Range("H77").Value = 753.75
$ws_ARM.Range("J77").Value = 0
$ws_ARM.Range("L77").Value = 0
$ws_ARM.Range("N77").ClearContents()
$ws_ARM.Range("H122").Value = 2856.1428
$ws_ARM.Range("I122").Value = 1499.25
$ws_ARM.Range("K122").Value = 4497.75
$ws_ARM.Range("M122").Value = -2047.75
$ws_ARM.Range("H132").Value = 1771.6666
$ws_ARM.Range("J132").Value = 2933.3333
$ws_ARM.Range("L132").Value = 8799.999899999999
$ws_ARM.Range("N132").Value = -13859.9999

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H41").Value = 2899.2856
$ws_CRP.Range("I41").Value = 59
$ws_CRP.Range("J41").Value = 10000
$ws_CRP.Range("K41").Value = 59
$ws_CRP.Range("L41").Value = 10000
$ws_CRP.Range("M41").Value = 369
$ws_CRP.Range("N41").Value = -10856
$ws_CRP.Range("H50").Value = 20082.75
$ws_CRP.Range("H59").Value = 28687.875
$ws_CRP.Range("I59").Value = 19752
$ws_CRP.Range("J59").Value = 31666.5
$ws_CRP.Range("K59").Value = 19752
$ws_CRP.Range("L59").Value = 31666.5
$ws_CRP.Range("M59").Value = -18607
$ws_CRP.Range("N59").Value = -33956.5
$ws_CRP.Range("H60").Value = 20181.666
$ws_CRP.Range("I60").Value = 20218
$ws_CRP.Range("K60").Value = 20218
$ws_CRP.Range("M60").Value = -19707

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H29").Value = 225.2
$ws_CUL.Range("I29").Value = 175.33333
$ws_CUL.Range("J29").Value = 300
$ws_CUL.Range("K29").Value = 525.99999
$ws_CUL.Range("L29").Value = 900
$ws_CUL.Range("M29").Value = -248.99999
$ws_CUL.Range("N29").Value = -1454
$ws_CUL.Range("H44").Value = 504.2
$ws_CUL.Range("I44").Value = 442.75
$ws_CUL.Range("J44").Value = 750
$ws_CUL.Range("K44").Value = 1328.25
$ws_CUL.Range("L44").Value = 2250
$ws_CUL.Range("M44").Value = -930.25
$ws_CUL.Range("N44").Value = -3046
$ws_CUL.Range("H52").Value = 1750
$ws_CUL.Range("J52").Value = 1750
$ws_CUL.Range("L52").Value = 5250
$ws_CUL.Range("N52").Value = -5782
$ws_CUL.Range("H92").Value = 1008.8
$ws_CUL.Range("I92").Value = 799
$ws_CUL.Range("K92").Value = 2397
$ws_CUL.Range("M92").Value = -1149
$ws_CUL.Range("H134").Value = 147338.72
$ws_CUL.Range("I134").Value = 250334.5
$ws_CUL.Range("K134").Value = 751003.5
$ws_CUL.Range("M134").Value = -745933.5

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H7").Value = 3000
$ws_GSM.Range("J7").Value = 3000
$ws_GSM.Range("L7").Value = 3000
$ws_GSM.Range("N7").Value = -3224
$ws_GSM.Range("H8").Value = 3000
$ws_GSM.Range("J8").Value = 3000
$ws_GSM.Range("L8").Value = 3000
$ws_GSM.Range("N8").Value = -3278
$ws_GSM.Range("H14").Value = 22999.5
$ws_GSM.Range("I14").Value = 20000
$ws_GSM.Range("J14").Value = 25999
$ws_GSM.Range("K14").Value = 20000
$ws_GSM.Range("L14").Value = 25999
$ws_GSM.Range("M14").Value = -19832
$ws_GSM.Range("N14").Value = -26335
$ws_GSM.Range("H20").Value = 38094.855
$ws_GSM.Range("J20").Value = 38094.855
$ws_GSM.Range("L20").Value = 38094.855
$ws_GSM.Range("N20").Value = -38584.855
$ws_GSM.Range("H33").Value = 0
$ws_GSM.Range("J33").Value = 0
$ws_GSM.Range("L33").Value = 0
$ws_GSM.Range("N33").ClearContents()

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H17").Value = 10750
$ws_LTW.Range("J17").Value = 9500
$ws_LTW.Range("L17").Value = 9500
$ws_LTW.Range("N17").Value = -9840
$ws_LTW.Range("H18").Value = 9500
$ws_LTW.Range("I18").Value = 9500
$ws_LTW.Range("K18").Value = 9500
$ws_LTW.Range("M18").Value = -9328
$ws_LTW.Range("H64").Value = 0
$ws_LTW.Range("J64").Value = 0
$ws_LTW.Range("L64").Value = 0
$ws_LTW.Range("N64").ClearContents()
$ws_LTW.Range("H67").Value = 0
$ws_LTW.Range("J67").Value = 0
$ws_LTW.Range("L67").Value = 0
$ws_LTW.Range("N67").ClearContents()
$ws_LTW.Range("H122").Value = 3590
$ws_LTW.Range("I122").Value = 3460
$ws_LTW.Range("J122").Value = 4500
$ws_LTW.Range("K122").Value = 10380
$ws_LTW.Range("L122").Value = 13500
$ws_LTW.Range("M122").Value = -7930
$ws_LTW.Range("N122").Value = -18400
$ws_LTW.Range("H132").Value = 9007.143
$ws_LTW.Range("J132").Value = 7694.25
$ws_LTW.Range("L132").Value = 23082.75
$ws_LTW.Range("N132").Value = -28142.75

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H122").Value = 2127.9375
$ws_WVR.Range("I122").Value = 2097.5715
$ws_WVR.Range("J122").Value = 2340.5
$ws_WVR.Range("K122").Value = 6292.7145
$ws_WVR.Range("L122").Value = 7021.5
$ws_WVR.Range("M122").Value = -3842.7145
$ws_WVR.Range("N122").Value = -11921.5

